$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update contract/offer reference values
$ws.Range("A2").Value = "ME-510"
$ws.Range("C2").Value = "CT-249"

# Update the selected cell on the sheet
$ws.Range("C3").Select()
